$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Vásárlás"

# Rename the 5th column header and give the first data row (Farkas Anikó) a value in it
$ws.Range("E1").Value = "Adás-Vételi"
$ws.Range("E2").Value = 6

# Append a new row (21st customer) duplicating Farkas Anikó's details
$ws.Range("A2:D2").Copy($ws.Range("A22:D22"))
$ws.Range("A22").Value = 21
$ws.Range("E22").Value = 10

# Match the column width that Excel ends up with (closest reachable value to 10.25)
$ws.Columns("E").ColumnWidth = 9.3

# Leave the cursor where the author ended up - on the new row
$ws.Range("D22").Select()
